$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "cunnect Node & re connet Node",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "cut information Node &  connet Node",
    2
)
